$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = "num vals for 180 degs"
$ws.Range("H10").Value = "earth circumference in meters"
$ws.Range("H11").Value = "meters per degree"
$ws.Range("H13").Value = "meters per tick"
$ws.Range("G8").Value = "ticks per degree"
$ws.Range("H5").Value = "num degrees"
$ws.Range("H14").Value = "mm per tick"

$ws.Range("G6").Formula = "=POWER(2, 31)"
$ws.Range("H6").Value = 180

$ws.Range("H8").Formula = "=G6/H6"

$ws.Range("G10").Value = 40075000

$ws.Range("G11").Formula = "=G10/360"

$ws.Range("G13").Formula = "=G11/H8"

$ws.Range("G14").Formula = "=G13*1000"

$ws.Columns.Item(7).ColumnWidth = 11.1640625
$ws.Columns.Item(8).ColumnWidth = 12.1640625

$ws.Range("H15").Select()
